$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$x = $s.NotesPage
$x.Shapes.AddTextbox(1, 10,10,10,10)
